# Insert a new row at position 239 (shifts existing rows 239:331 down to 240:332)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("239:239").Insert()

# Populate the newly inserted row 239 with the new record
$ws.Range("A239").Value = 4
$ws.Range("B239").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C239").Value = "Los Lagos"
$ws.Range("D239").Value = 44795
$ws.Range("E239").Value = 10
$ws.Range("F239").Value = 100112040
$ws.Range("G239").Value = "Cilantro"
$ws.Range("H239").Value = "Sin especificar"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 70
$ws.Range("K239").Value = 11500
$ws.Range("L239").Value = 11500
$ws.Range("M239").Value = 11500
$ws.Range("N239").Value = '$/caja 36 atados'
$ws.Range("O239").Value = "Región Metropolitana"
$ws.Range("P239").Value = 319
$ws.Range("Q239").Value = 36
$ws.Range("R239").Value = "Hortaliza"
